$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 39-41 ---
# Row 39 (WE 12/5/20): Positive count 385 -> 390
$ws.Range("B39").Value = 390

# Row 40 (WE 12/12/20): Positive count 406 -> 411
$ws.Range("B40").Value = 411

# Row 41 (WE 12/19/20): Positive count 210 -> 407 (verified/corrected),
# and it's no longer the latest partial-week row, so drop the red
# "preliminary" font formatting back to the normal date style.
$ws.Range("B41").Value = 407
$ws.Range("A41").Style = "Normal"

# --- Append new rows for newly reported weeks ---
# Row 42 (WE 12/26/20)
$ws.Range("A42").Value = 44191
$ws.Range("A42").NumberFormat = "m/d/yy"
$ws.Range("B42").Value = 419
$ws.Range("C42").Formula = "=(B42-B41)/B42"
$ws.Range("D42").Formula = "=(B42/79546)*100000"

# Row 43 (WE 1/2/21)
$ws.Range("A43").Value = 44198
$ws.Range("A43").NumberFormat = "m/d/yy"
$ws.Range("B43").Value = 443
$ws.Range("C43").Formula = "=(B43-B42)/B43"
$ws.Range("D43").Formula = "=(B43/79546)*100000"

# Row 44 (WE 1/9/21) - partial week, newest data point -> highlight red
$ws.Range("A44").Value = 44205
$ws.Range("A44").NumberFormat = "m/d/yy"
$ws.Range("B44").Value = 405
$ws.Range("C44").Formula = "=(B44-B43)/B44"
$ws.Range("D44").Formula = "=(B44/79546)*100000"

$ws.Range("A44").Font.Color = 255
$ws.Range("A44").Interior.ColorIndex = $ws.Range("A41").Interior.ColorIndex

# Copy the cell style from rows 39-40 (normal, non-highlighted date format)
# onto the new, non-final rows and keep highlighting only on the newest row.
$ws.Range("A39:A40").Copy() | Out-Null
$ws.Range("A42:A43").PasteSpecial(-4122) | Out-Null

$ws.Range("A39").Copy() | Out-Null
$ws.Range("A41").PasteSpecial(-4122) | Out-Null

$ws.Range("A41").Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null

$ws.Range("C39:D40").Copy() | Out-Null
$ws.Range("C42:D44").PasteSpecial(-4122) | Out-Null

# --- Update the sheet view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("B40").Select() | Out-Null

$wb.Save()
